# Apply postgame hitter report corrections (Lavoie, Brady - 2022-07-29)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Pitch block @ row 10-17 (Pitcher: row 14 name) ---
$ws.Range("J10").Value = 2
$ws.Range("M10").Value = ""
$ws.Range("M12").Value = ""
$ws.Range("J14").Value = "Roblez"
$ws.Range("M14").Value = ""
$ws.Range("J16").Value = "88-90 MPH"
$ws.Range("J17").Value = "FB,CB,CH"

# --- Pitch block @ row 19-26 ---
$ws.Range("M19").Value = ""
$ws.Range("M21").Value = ""
$ws.Range("J26").Value = "FB,CB,CH"

# --- Pitch block @ row 28-35 (Pitcher: row 32 name) ---
$ws.Range("J28").Value = 4
$ws.Range("M28").Value = ""
$ws.Range("J29").Value = 2
$ws.Range("M30").Value = ""
$ws.Range("J32").Value = "Herbst"
$ws.Range("M32").Value = "Fly Ball"
$ws.Range("M33").Value = "Out"
$ws.Range("J34").Value = "83-85 MPH"
$ws.Range("J35").Value = "SL,FB,CB,CH"

# --- Pitch block @ row 37-44 (Pitcher: row 41 name) ---
$ws.Range("J37").Value = 6
$ws.Range("M37").Value = ""
$ws.Range("J38").Value = 0
$ws.Range("M39").Value = ""
$ws.Range("J41").Value = "Herbst"
$ws.Range("M41").Value = ""
$ws.Range("J43").Value = "83-85 MPH"
$ws.Range("J44").Value = "SL,FB,CB,CH"

# --- Pitch block @ row 46-53 (Pitcher: row 50 name) ---
$ws.Range("J46").Value = 7
$ws.Range("M46").Value = ""
$ws.Range("J47").Value = 1
$ws.Range("M48").Value = ""
$ws.Range("J50").Value = "Plum"
$ws.Range("M50").Value = ""
$ws.Range("M51").Value = "Undefined"
$ws.Range("J52").Value = "84-86 MPH"
$ws.Range("J53").Value = "SL,FB,CH"

# --- Pitch block @ row 61-65 ---
$ws.Range("M61").Value = ""
$ws.Range("M63").Value = ""
$ws.Range("M65").Value = ""
